# Applies the edit described by the diff:
#   Two new price-report rows are inserted right before row 473 (pushing the
#   existing rows 473-568 down to 475-570), and the two freshly inserted
#   rows (473-474) are populated with a new "Limón" price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 473 - this shifts the former rows
# 473:568 down to 475:570, matching the dimension change (T568 -> T570).
$ws.Rows("473:474").Insert()

# --- New row 473 ---------------------------------------------------------
$ws.Cells.Item(473, 1).Value  = 4
$ws.Cells.Item(473, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(473, 3).Value  = "Los Lagos"
$ws.Cells.Item(473, 4).Value  = 44785
$ws.Cells.Item(473, 5).Value  = 10
$ws.Cells.Item(473, 6).Value  = "Fruta"
$ws.Cells.Item(473, 7).Value  = 100102
$ws.Cells.Item(473, 8).Value  = "Cítricos"
$ws.Cells.Item(473, 9).Value  = 100102003
$ws.Cells.Item(473, 10).Value = "Limón"
$ws.Cells.Item(473, 11).Value = "Sin especificar"
$ws.Cells.Item(473, 12).Value = "1a amarillo"
$ws.Cells.Item(473, 13).Value = 1400
$ws.Cells.Item(473, 14).Value = 8000
$ws.Cells.Item(473, 15).Value = 9000
$ws.Cells.Item(473, 16).Value = 8500
$ws.Cells.Item(473, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(473, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(473, 19).Value = 472
$ws.Cells.Item(473, 20).Value = 18

# --- New row 474 ---------------------------------------------------------
$ws.Cells.Item(474, 1).Value  = 4
$ws.Cells.Item(474, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(474, 3).Value  = "Los Lagos"
$ws.Cells.Item(474, 4).Value  = 44785
$ws.Cells.Item(474, 5).Value  = 10
$ws.Cells.Item(474, 6).Value  = "Fruta"
$ws.Cells.Item(474, 7).Value  = 100102
$ws.Cells.Item(474, 8).Value  = "Cítricos"
$ws.Cells.Item(474, 9).Value  = 100102003
$ws.Cells.Item(474, 10).Value = "Limón"
$ws.Cells.Item(474, 11).Value = "Sin especificar"
$ws.Cells.Item(474, 12).Value = "2a amarillo"
$ws.Cells.Item(474, 13).Value = 700
$ws.Cells.Item(474, 14).Value = 6000
$ws.Cells.Item(474, 15).Value = 6000
$ws.Cells.Item(474, 16).Value = 6000
$ws.Cells.Item(474, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(474, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(474, 19).Value = 333
$ws.Cells.Item(474, 20).Value = 18
